$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "cotinine_level"
$ws.Range("B1").Value = "group"
$ws.Range("C1").Value = "description"

$ws.Range("C1").Select()
